# Replace all occurrences of "BMP MONEY PLUS SOCIEDADE DE CRÉDITO DIRETO S.A"
# with "BMP SOCIEDADE DE CRÉDITO DIRETO S.A" throughout the document body,
# per commit message:
#   "alteracao BMP MONEY PLUS SOCIEDADE DE CRÉDITO DIRETO S.A para
#    BMP SOCIEDADE DE CRÉDITO DIRETO S.A"

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# wdFindContinue = 1, wdReplaceAll = 2
$find.Execute(
    "BMP MONEY PLUS SOCIEDADE DE CRÉDITO DIRETO S.A",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "BMP SOCIEDADE DE CRÉDITO DIRETO S.A",
    2
)
